$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.503.52'
$ws.Range("E2").Value = '  +3.33%  '
$ws.Range("D3").Value = '3.076.59'
$ws.Range("E3").Value = '  +5.54%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.62%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +3.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.29'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.85%  '
$ws.Range("E10").Value = '  +4.13%  '
$ws.Range("E11").Value = '  +6.65%  '
$ws.Range("D12").Value = '3.606.01'
$ws.Range("E12").Value = '  +5.37%  '
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000165'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.14%  '
$ws.Range("D16").Value = '57.596.38'
$ws.Range("E16").Value = '  +3.44%  '
$ws.Range("D17").Value = '3.080.48'
$ws.Range("E17").Value = '  +5.44%  '
$ws.Range("E18").Value = '  +2.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.86%  '
$ws.Range("E20").Value = '  +6.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '338.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.16%  '
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.500'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.33%  '
$ws.Range("E25").Value = '  +7.07%  '
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").Value = '0.0₃0943'
$ws.Range("E27").Value = '  +12.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.48%  '
$ws.Range("E30").Value = '  +3.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.89'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0678'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.26%  '
$ws.Range("D39").Value = '3.114.31'
$ws.Range("E39").Value = '  +5.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.672'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.35%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.47%  '
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").Value = '2.263.85'
$ws.Range("E44").Value = '  +6.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.82%  '
$ws.Range("E46").Value = '  +6.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.956'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.56%  '
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("E50").Value = '  +4.38%  '
$ws.Range("E51").Value = '  +4.23%  '
